# Insert a new row at 171, shifting existing rows 171-271 down to 172-272
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(171).Insert()

# Populate the newly inserted row 171 with the new record's data
$ws.Range("A171").Value = 5
$ws.Range("B171").Value = "Macroferia Regional de Talca"
$ws.Range("C171").Value = "Maule"
$ws.Range("D171").Value = "2022-02-14"
$ws.Range("E171").Value = 7
$ws.Range("F171").Value = 100112023
$ws.Range("G171").Value = "Brócoli"
$ws.Range("H171").Value = "Sin especificar"
$ws.Range("I171").Value = "Primera"
$ws.Range("J171").Value = 3000
$ws.Range("K171").Value = 800
$ws.Range("L171").Value = 800
$ws.Range("M171").Value = 800
$ws.Range("N171").Value = "$/unidad"
$ws.Range("O171").Value = "Región del Maule"
$ws.Range("P171").Value = 800
$ws.Range("Q171").Value = 1
$ws.Range("R171").Value = "Hortaliza"
